$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data per daily GitHub Actions refresh.
# Also swap rows 48/49 (EnergySwap and Aptos changed rank order).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.204.60'
$ws.Range("E2").Value = '  +0.13%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.870.82'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.69'
$ws.Range("E5").Value = '  -0.09%  '

$ws.Range("E6").Value = '  +0.11%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4698'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2844'
$ws.Range("E8").Value = '  -1.48%  '

$ws.Range("E9").Value = '  +0.07%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.22'
$ws.Range("E10").Value = '  -2.77%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07786'
$ws.Range("E11").Value = '  -2.00%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '96.29'
$ws.Range("E12").Value = '  -1.14%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.872.92'
$ws.Range("E13").Value = '  +1.26%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6896'
$ws.Range("E14").Value = '  +2.29%  '

$ws.Range("E15").Value = '  -0.02%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '266.52'
$ws.Range("E16").Value = '  -0.61%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.192.62'
$ws.Range("E17").Value = '  +0.20%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.71'
$ws.Range("E18").Value = '  +0.76%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007700'
$ws.Range("E19").Value = '  +1.00%  '

$ws.Range("E20").Value = '  +0.09%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.114.74'
$ws.Range("E21").Value = '  +0.83%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.002'
$ws.Range("E22").Value = '  +0.17%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.233'
$ws.Range("E23").Value = '  +1.02%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.159'
$ws.Range("E24").Value = '  +0.54%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.490'
$ws.Range("E25").Value = '  +3.85%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.34'
$ws.Range("E26").Value = '  -0.51%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.73'
$ws.Range("E27").Value = '  -0.23%  '

$ws.Range("E28").Value = '  +0.58%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.371'
$ws.Range("E29").Value = '  -0.54%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09902'
$ws.Range("E30").Value = '  +0.91%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.354'
$ws.Range("E31").Value = '  +1.90%  '

$ws.Range("E32").Value = '  -0.24%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.045'
$ws.Range("E33").Value = '  +1.52%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04729'
$ws.Range("E34").Value = '  +0.83%  '

$ws.Range("E35").Value = '  +1.33%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6992'
$ws.Range("E36").Value = '  +0.33%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.719'
$ws.Range("E37").Value = '  +0.52%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01861'
$ws.Range("E38").Value = '  -0.04%  '

$ws.Range("E39").Value = '  +7.17%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.242'
$ws.Range("E40").Value = '  -1.12%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '72.61'
$ws.Range("E41").Value = '  -0.28%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.935'
$ws.Range("E42").Value = '  +0.36%  '

$ws.Range("E43").Value = '  +0.21%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4144'
$ws.Range("E44").Value = '  +0.57%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8318'
$ws.Range("E45").Value = '  -0.53%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '102.86'
$ws.Range("E46").Value = '  -0.04%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '973.92'
$ws.Range("E47").Value = '  +3.75%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.176'
$ws.Range("E48").Value = '  -0.03%  '

$ws.Range("B49").Value = 'Aptos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.053'
$ws.Range("E49").Value = '  +0.87%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.44'
$ws.Range("E50").Value = '  +1.93%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05646'
$ws.Range("E51").Value = '  -0.14%  '

